$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the "Lily" translation into column D for row 3
$ws.Range("D3").Value = "Lily"

# Add new translation "Lime" for the enemy/character in row 26 (column D)
$ws.Range("D26").Value = "Lime"

# Duplicate the "Shina" translation into column D for row 30
$ws.Range("D30").Value = "Shina"
